$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 10; existing rows 10.. shift down to 11..,
# pushing the former last row (58) down to the new last row (59).
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new weekly price record.
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C10").Value = "Arica y Parinacota"
$ws.Range("D10").Value = 44547
$ws.Range("E10").Value = 15
$ws.Range("F10").Value = 100112038
$ws.Range("G10").Value = "Cebollín baby"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 300
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 3000
$ws.Range("M10").Value = 2750
$ws.Range("N10").Value = "$/paquete 1,5 a 2 kilos"
$ws.Range("O10").Value = "Región de Arica y Parinacota"
$ws.Range("P10").Value = 1375
$ws.Range("Q10").Value = 2
$ws.Range("R10").Value = "Hortaliza"
